$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("B9").Value = 6865281
$ws.Range("C9").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D9").Value = 45150.5
$ws.Range("E9").Value = 'GOSK Gabela'
$ws.Range("F9").Value = 'Zvijezda 09'
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 'H'
$ws.Range("J9").Value = 1.75
$ws.Range("K9").Value = 4
$ws.Range("L9").Value = 3.5
$ws.Range("M9").Value = 1.75
$ws.Range("N9").Value = 4
$ws.Range("O9").Value = 3.4
$ws.Range("P9").Value = -0.5
$ws.Range("Q9").Value = 1.8
$ws.Range("R9").Value = 2
$ws.Range("S9").Value = 2.5
$ws.Range("T9").Value = 1.85
$ws.Range("U9").Value = 1.95
$ws.Range("V9").Value = 0.75
$ws.Range("W9").Value = -1
$ws.Range("X9").Value = -1
$ws.Range("Y9").Value = 0.8
$ws.Range("Z9").Value = -1
$ws.Range("AA9").Value = -1
$ws.Range("AB9").Value = 0.95

# Row 10
$ws.Range("B10").Value = 6865285
$ws.Range("C10").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D10").Value = 45150.5
$ws.Range("E10").Value = 'NK Igman Konjic'
$ws.Range("F10").Value = 'Sloga'
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 'H'
$ws.Range("J10").Value = 2
$ws.Range("K10").Value = 3.4
$ws.Range("L10").Value = 3.2
$ws.Range("M10").Value = 1.909
$ws.Range("N10").Value = 3.5
$ws.Range("O10").Value = 3.4
$ws.Range("P10").Value = -0.5
$ws.Range("Q10").Value = 1.95
$ws.Range("R10").Value = 1.85
$ws.Range("S10").Value = 2.5
$ws.Range("T10").Value = 1.85
$ws.Range("U10").Value = 1.95
$ws.Range("V10").Value = 0.909
$ws.Range("W10").Value = -1
$ws.Range("X10").Value = -1
$ws.Range("Y10").Value = 0.95
$ws.Range("Z10").Value = -1
$ws.Range("AA10").Value = -1
$ws.Range("AB10").Value = 0.95

# Row 29
$ws.Range("B29").Value = 6865295
$ws.Range("C29").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D29").Value = 45172.61458333334
$ws.Range("E29").Value = 'FK Tuzla City'
$ws.Range("F29").Value = 'NK Igman Konjic'
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 1
$ws.Range("I29").Value = 'H'
$ws.Range("J29").Value = 1.8
$ws.Range("K29").Value = 3.4
$ws.Range("L29").Value = 3.8
$ws.Range("M29").Value = 1.615
$ws.Range("N29").Value = 3.5
$ws.Range("O29").Value = 4.5
$ws.Range("P29").Value = -0.75
$ws.Range("Q29").Value = 1.85
$ws.Range("R29").Value = 1.95
$ws.Range("S29").Value = 2.75
$ws.Range("T29").Value = 2
$ws.Range("U29").Value = 1.8
$ws.Range("V29").Value = 0.615
$ws.Range("W29").Value = -1
$ws.Range("X29").Value = -1
$ws.Range("Y29").Value = 0.8500000000000001
$ws.Range("Z29").Value = -1
$ws.Range("AA29").Value = 1
$ws.Range("AB29").Value = -1

# Row 30
$ws.Range("B30").Value = 6865296
$ws.Range("C30").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D30").Value = 45172.61458333334
$ws.Range("E30").Value = 'Velez Mostar'
$ws.Range("F30").Value = 'Zeljeznicar'
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 'H'
$ws.Range("J30").Value = 1.909
$ws.Range("K30").Value = 3.2
$ws.Range("L30").Value = 3.6
$ws.Range("M30").Value = 1.95
$ws.Range("N30").Value = 3.2
$ws.Range("O30").Value = 3.4
$ws.Range("P30").Value = -0.5
$ws.Range("Q30").Value = 2.025
$ws.Range("R30").Value = 1.775
$ws.Range("S30").Value = 2.25
$ws.Range("T30").Value = 1.9
$ws.Range("U30").Value = 1.9
$ws.Range("V30").Value = 0.95
$ws.Range("W30").Value = -1
$ws.Range("X30").Value = -1
$ws.Range("Y30").Value = 1.025
$ws.Range("Z30").Value = -1
$ws.Range("AA30").Value = -1
$ws.Range("AB30").Value = 0.8999999999999999

# Row 49
$ws.Range("B49").Value = 6865310
$ws.Range("C49").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D49").Value = 45200.41666666666
$ws.Range("E49").Value = 'NK Igman Konjic'
$ws.Range("F49").Value = 'Zrinjski Mostar'
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 2
$ws.Range("I49").Value = 'A'
$ws.Range("J49").Value = 3.4
$ws.Range("K49").Value = 3.6
$ws.Range("L49").Value = 1.833
$ws.Range("M49").Value = 4.75
$ws.Range("N49").Value = 4.75
$ws.Range("O49").Value = 1.45
$ws.Range("P49").Value = 1.25
$ws.Range("Q49").Value = 1.775
$ws.Range("R49").Value = 2.025
$ws.Range("S49").Value = 2.75
$ws.Range("T49").Value = 1.85
$ws.Range("U49").Value = 1.95
$ws.Range("V49").Value = -1
$ws.Range("W49").Value = -1
$ws.Range("X49").Value = 0.45
$ws.Range("Y49").Value = -1
$ws.Range("Z49").Value = 1.025
$ws.Range("AA49").Value = -1
$ws.Range("AB49").Value = 0.95

# Row 50
$ws.Range("B50").Value = 6865311
$ws.Range("C50").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D50").Value = 45200.41666666666
$ws.Range("E50").Value = 'Sloga'
$ws.Range("F50").Value = 'GOSK Gabela'
$ws.Range("G50").Value = 3
$ws.Range("H50").Value = 2
$ws.Range("I50").Value = 'H'
$ws.Range("J50").Value = 1.833
$ws.Range("K50").Value = 3.6
$ws.Range("L50").Value = 3.4
$ws.Range("M50").Value = 1.909
$ws.Range("N50").Value = 3.4
$ws.Range("O50").Value = 3.3
$ws.Range("P50").Value = -0.5
$ws.Range("Q50").Value = 1.925
$ws.Range("R50").Value = 1.875
$ws.Range("S50").Value = 2.25
$ws.Range("T50").Value = 1.825
$ws.Range("U50").Value = 1.975
$ws.Range("V50").Value = 0.909
$ws.Range("W50").Value = -1
$ws.Range("X50").Value = -1
$ws.Range("Y50").Value = 0.925
$ws.Range("Z50").Value = -1
$ws.Range("AA50").Value = 0.825
$ws.Range("AB50").Value = -1

# Row 76
$ws.Range("B76").Value = 6865377
$ws.Range("C76").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D76").Value = 45235.375
$ws.Range("E76").Value = 'Zrinjski Mostar'
$ws.Range("F76").Value = 'FK Tuzla City'
$ws.Range("G76").Value = 3
$ws.Range("H76").Value = 1
$ws.Range("I76").Value = 'H'
$ws.Range("J76").Value = 1.333
$ws.Range("K76").Value = 5
$ws.Range("L76").Value = 6
$ws.Range("M76").Value = 1.166
$ws.Range("N76").Value = 6.5
$ws.Range("O76").Value = 13
$ws.Range("P76").Value = -2
$ws.Range("Q76").Value = 1.9
$ws.Range("R76").Value = 1.9
$ws.Range("S76").Value = 3.25
$ws.Range("T76").Value = 1.95
$ws.Range("U76").Value = 1.85
$ws.Range("V76").Value = 0.1659999999999999
$ws.Range("W76").Value = -1
$ws.Range("X76").Value = -1
$ws.Range("Y76").Value = 0
$ws.Range("Z76").Value = 0
$ws.Range("AA76").Value = 0.95
$ws.Range("AB76").Value = -1

# Row 77
$ws.Range("B77").Value = 6865328
$ws.Range("C77").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D77").Value = 45235.375
$ws.Range("E77").Value = 'Siroki Brijeg'
$ws.Range("F77").Value = 'NK Posusje'
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = 1
$ws.Range("I77").Value = 'D'
$ws.Range("J77").Value = 2
$ws.Range("K77").Value = 3
$ws.Range("L77").Value = 3.5
$ws.Range("M77").Value = 2.1
$ws.Range("N77").Value = 3
$ws.Range("O77").Value = 3.3
$ws.Range("P77").Value = -0.25
$ws.Range("Q77").Value = 1.825
$ws.Range("R77").Value = 1.975
$ws.Range("S77").Value = 2
$ws.Range("T77").Value = 1.825
$ws.Range("U77").Value = 1.975
$ws.Range("V77").Value = -1
$ws.Range("W77").Value = 2
$ws.Range("X77").Value = -1
$ws.Range("Y77").Value = -0.5
$ws.Range("Z77").Value = 0.4875
$ws.Range("AA77").Value = 0
$ws.Range("AB77").Value = 0

# Row 122
$ws.Range("B122").Value = 6865381
$ws.Range("C122").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D122").Value = 45353.375
$ws.Range("E122").Value = 'FK Tuzla City'
$ws.Range("F122").Value = 'Zvijezda 09'
$ws.Range("G122").Value = 2
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 'H'
$ws.Range("J122").Value = 1.666
$ws.Range("K122").Value = 3.6
$ws.Range("L122").Value = 4.333
$ws.Range("M122").Value = 1.5
$ws.Range("N122").Value = 4
$ws.Range("O122").Value = 5.25
$ws.Range("P122").Value = -1
$ws.Range("Q122").Value = 1.925
$ws.Range("R122").Value = 1.875
$ws.Range("S122").Value = 2.5
$ws.Range("T122").Value = 1.8
$ws.Range("U122").Value = 2
$ws.Range("V122").Value = 0.5
$ws.Range("W122").Value = -1
$ws.Range("X122").Value = -1
$ws.Range("Y122").Value = 0.925
$ws.Range("Z122").Value = -1
$ws.Range("AA122").Value = -1
$ws.Range("AB122").Value = 1

# Row 123
$ws.Range("B123").Value = 6865363
$ws.Range("C123").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D123").Value = 45353.375
$ws.Range("E123").Value = 'NK Igman Konjic'
$ws.Range("F123").Value = 'Siroki Brijeg'
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 'H'
$ws.Range("J123").Value = 2
$ws.Range("K123").Value = 3.3
$ws.Range("L123").Value = 3.25
$ws.Range("M123").Value = 2.3
$ws.Range("N123").Value = 3.2
$ws.Range("O123").Value = 2.75
$ws.Range("P123").Value = -0.25
$ws.Range("Q123").Value = 2.05
$ws.Range("R123").Value = 1.75
$ws.Range("S123").Value = 2
$ws.Range("T123").Value = 1.9
$ws.Range("U123").Value = 1.9
$ws.Range("V123").Value = 1.3
$ws.Range("W123").Value = -1
$ws.Range("X123").Value = -1
$ws.Range("Y123").Value = 1.05
$ws.Range("Z123").Value = -1
$ws.Range("AA123").Value = -1
$ws.Range("AB123").Value = 0.8999999999999999

# Row 189
$ws.Range("B189").Value = 7952780
$ws.Range("C189").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D189").Value = 45432.5
$ws.Range("E189").Value = 'Velez Mostar'
$ws.Range("F189").Value = 'GOSK Gabela'
$ws.Range("G189").Value = 3
$ws.Range("H189").Value = 3
$ws.Range("I189").Value = 'D'
$ws.Range("J189").Value = 1.4
$ws.Range("K189").Value = 4
$ws.Range("L189").Value = 7
$ws.Range("M189").Value = 1.363
$ws.Range("N189").Value = 4.2
$ws.Range("O189").Value = 8
$ws.Range("P189").Value = -1.5
$ws.Range("Q189").Value = 2
$ws.Range("R189").Value = 1.8
$ws.Range("S189").Value = 2.75
$ws.Range("T189").Value = 1.825
$ws.Range("U189").Value = 1.975
$ws.Range("V189").Value = -1
$ws.Range("W189").Value = 3.2
$ws.Range("X189").Value = -1
$ws.Range("Y189").Value = -1
$ws.Range("Z189").Value = 0.8
$ws.Range("AA189").Value = 0.825
$ws.Range("AB189").Value = -1

# Row 191
$ws.Range("B191").Value = 7952781
$ws.Range("C191").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D191").Value = 45432.5
$ws.Range("E191").Value = 'Zvijezda 09'
$ws.Range("F191").Value = 'Zeljeznicar'
$ws.Range("G191").Value = 0
$ws.Range("H191").Value = 5
$ws.Range("I191").Value = 'A'
$ws.Range("J191").Value = 2.15
$ws.Range("K191").Value = 3.25
$ws.Range("L191").Value = 2.9
$ws.Range("M191").Value = 3.6
$ws.Range("N191").Value = 3.4
$ws.Range("O191").Value = 1.85
$ws.Range("P191").Value = 0.5
$ws.Range("Q191").Value = 1.875
$ws.Range("R191").Value = 1.925
$ws.Range("S191").Value = 2.5
$ws.Range("T191").Value = 1.975
$ws.Range("U191").Value = 1.825
$ws.Range("V191").Value = -1
$ws.Range("W191").Value = -1
$ws.Range("X191").Value = 0.8500000000000001
$ws.Range("Y191").Value = -1
$ws.Range("Z191").Value = 0.925
$ws.Range("AA191").Value = 0.9750000000000001
$ws.Range("AB191").Value = -1

# Row 193
$ws.Range("B193").Value = 7952778
$ws.Range("C193").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D193").Value = 45432.5
$ws.Range("E193").Value = 'Sloga'
$ws.Range("F193").Value = 'Siroki Brijeg'
$ws.Range("G193").Value = 2
$ws.Range("H193").Value = 3
$ws.Range("I193").Value = 'A'
$ws.Range("J193").Value = 1.727
$ws.Range("K193").Value = 3.75
$ws.Range("L193").Value = 3.75
$ws.Range("M193").Value = 1.7
$ws.Range("N193").Value = 3.9
$ws.Range("O193").Value = 3.9
$ws.Range("P193").Value = -0.75
$ws.Range("Q193").Value = 1.975
$ws.Range("R193").Value = 1.825
$ws.Range("S193").Value = 2.25
$ws.Range("T193").Value = 1.8
$ws.Range("U193").Value = 2
$ws.Range("V193").Value = -1
$ws.Range("W193").Value = -1
$ws.Range("X193").Value = 2.9
$ws.Range("Y193").Value = -1
$ws.Range("Z193").Value = 0.825
$ws.Range("AA193").Value = 0.8
$ws.Range("AB193").Value = -1

# Row 194
$ws.Range("B194").Value = 7952777
$ws.Range("C194").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D194").Value = 45432.5
$ws.Range("E194").Value = 'Borac Banja Luka'
$ws.Range("F194").Value = 'NK Igman Konjic'
$ws.Range("G194").Value = 4
$ws.Range("H194").Value = 3
$ws.Range("I194").Value = 'H'
$ws.Range("J194").Value = 1.25
$ws.Range("K194").Value = 5.75
$ws.Range("L194").Value = 7
$ws.Range("M194").Value = 1.2
$ws.Range("N194").Value = 5.75
$ws.Range("O194").Value = 12
$ws.Range("P194").Value = -2
$ws.Range("Q194").Value = 1.95
$ws.Range("R194").Value = 1.85
$ws.Range("S194").Value = 3.25
$ws.Range("T194").Value = 1.9
$ws.Range("U194").Value = 1.9
$ws.Range("V194").Value = 0.2
$ws.Range("W194").Value = -1
$ws.Range("X194").Value = -1
$ws.Range("Y194").Value = -1
$ws.Range("Z194").Value = 0.8500000000000001
$ws.Range("AA194").Value = 0.8999999999999999
$ws.Range("AB194").Value = -1
